$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Cd9","L1cam","ECs",3,1,146.0459156666667,438.137747,0.4862506770104965,0.4862506770104965,3,1,19.72083766666667,59.162513,0.8016210077351786,0.8016210077351787,2880.147794742024,25921.33015267821,0.389788757717067,0.3897887577170671),
    @("ECs","Cd9","L1cam","FAPs",3,1,146.0459156666667,438.137747,0.4862506770104965,0.4862506770104965,1,0.3333333333333333,0.099159,0.297477,0.004030657259573097,0.004030657259573097,14.481766951591,130.335902564319,0.00195990982126469,0.001959909821264691),
    @("ECs","Cd9","L1cam","sCs",3,1,146.0459156666667,438.137747,0.4862506770104965,0.4862506770104965,3,1,4.781202,14.343606,0.1943483350052483,0.1943483350052483,698.2750240772981,6284.475216695682,0.09450200947216474,0.09450200947216475),
    @("FAPs","Cd9","L1cam","ECs",3,1,122.7232436666666,368.169731,0.408599309644787,0.408599309644787,3,1,19.72083766666667,59.162513,0.8016210077351786,0.8016210077351787,2420.205166277111,21781.846496494,0.3275417903573524,0.3275417903573525),
    @("FAPs","Cd9","L1cam","FAPs",3,1,122.7232436666666,368.169731,0.408599309644787,0.408599309644787,1,0.3333333333333333,0.099159,0.297477,0.004030657259573097,0.004030657259573097,12.169114118743,109.522027068687,0.001646923773676316,0.001646923773676317),
    @("FAPs","Cd9","L1cam","sCs",3,1,122.7232436666666,368.169731,0.408599309644787,0.408599309644787,3,1,4.781202,14.343606,0.1943483350052483,0.1943483350052483,586.7646180655539,5280.881562589986,0.07941059551375822,0.07941059551375824),
    @("sCs","Cd9","L1cam","ECs",3,1,31.58192,94.74576,0.1051500133447165,0.1051500133447165,3,1,19.72083766666667,59.162513,0.8016210077351786,0.8016210077351787,622.8219175216534,5605.397257694881,0.0842904596607591,0.08429045966075913),
    @("sCs","Cd9","L1cam","FAPs",3,1,31.58192,94.74576,0.1051500133447165,0.1051500133447165,1,0.3333333333333333,0.099159,0.297477,0.004030657259573097,0.004030657259573097,3.13163160528,28.18468444752,0.0004238236646320895,0.0004238236646320896),
    @("sCs","Cd9","L1cam","sCs",3,1,31.58192,94.74576,0.1051500133447165,0.1051500133447165,3,1,4.781202,14.343606,0.1943483350052483,0.1943483350052483,150.99953906784,1358.99585161056,0.02043573001932528,0.02043573001932529)
)

$nRows = $data.Count
$nCols = 20
$arr = New-Object "object[,]" $nRows,$nCols
for ($r = 0; $r -lt $nRows; $r++) {
    for ($c = 0; $c -lt $nCols; $c++) {
        $arr[$r,$c] = $data[$r][$c]
    }
}

$startRow = 2
$endRow = $startRow + $nRows - 1
$rng = $ws.Range("A2:T" + $endRow)
$rng.Value = $arr

Write-Output "applied"
